$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("C8").Value = 45994
$ws.Range("D8").Value = 85000

# Row 9
$ws.Range("B9").Value = "DARWIN FUTBOL"
$ws.Range("C9").Value = 45921
$ws.Range("D9").Value = 200000

# Row 10
$ws.Range("B10").Value = "DAVIDCITO"
$ws.Range("C10").Value = 45947
$ws.Range("D10").Value = 100000

# Row 11
$ws.Range("B11").Value = "EL RUBY"
$ws.Range("C11").Value = 45992
$ws.Range("D11").Value = 85100

# Row 12
$ws.Range("B12").Value = "LA PAMPA"
$ws.Range("C12").Value = 45994
$ws.Range("D12").Value = 249000

# Row 13
$ws.Range("B13").Value = "LA SELECTA"
$ws.Range("C13").Value = 45912
$ws.Range("D13").Value = 82000

# Row 14
$ws.Range("B14").Value = "MARIANA"
$ws.Range("C14").Value = 45650
$ws.Range("D14").Value = 171900

# Row 15
$ws.Range("B15").Value = "MERKA FRUVER ALEJANDRO"
$ws.Range("C15").Value = 45988
$ws.Range("D15").Value = 60900

# Row 16
$ws.Range("B16").Value = "MERKA FRUVER MILDRED"
$ws.Range("D16").Value = 115400

# Row 17
$ws.Range("B17").Value = "MEZA 2"
$ws.Range("C17").Value = 45989
$ws.Range("D17").Value = 188000

# Row 18
$ws.Range("B18").Value = "MULTICARNES"
$ws.Range("D18").Value = 558300

# Row 19
$ws.Range("B19").Value = "NOVILLON SAN MATEO"
$ws.Range("C19").Value = 45971
$ws.Range("D19").Value = 83000

# Row 20
$ws.Range("B20").Value = "PINILLA"
$ws.Range("C20").Value = 45931
$ws.Range("D20").Value = 166000

# Row 21
$ws.Range("C21").Value = 45924
$ws.Range("D21").Value = 16000

# Row 22
$ws.Range("B22").Value = "PINILLA SOACHA"
$ws.Range("C22").Value = 45993
$ws.Range("D22").Value = 129000

# Row 23
$ws.Range("B23").Value = "PLAZA JESSICA"
$ws.Range("D23").Value = 621000

# Row 24
$ws.Range("C24").Value = 45995
$ws.Range("D24").Value = 1580300
